# "july 25 update 2" - add the two new days (rows 24 & 25) of Babine
# fence-count data and fill the running-total formulas down to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Babine")

# ---- Row 24 : 2024-07-25 (serial 45496) --------------------------------
$ws.Range("B24").Value = 25407
$ws.Range("C24").Value = 1060
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 29
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0

$ws.Range("L24").Formula = "=L23+B24"
$ws.Range("M24").Formula = "=I24+K24"
$ws.Range("N24").Formula = "=B24+M24"
$ws.Range("O24").Formula = "=O23+N24"
$ws.Range("P24").Formula = "=C24+J24"
$ws.Range("Q24").Formula = "=P24+Q23"

# ---- Row 25 : 2024-07-26 (serial 45497) --------------------------------
$ws.Range("B25").Value = 13068
$ws.Range("C25").Value = 538
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0

$ws.Range("L25").Formula = "=L24+B25"
$ws.Range("M25").Formula = "=I25+K25"
$ws.Range("N25").Formula = "=B25+M25"
$ws.Range("O25").Formula = "=O24+N25"
$ws.Range("P25").Formula = "=C25+J25"
$ws.Range("Q25").Formula = "=P25+Q24"

# ---- view / selection tidy-up (matches the author's re-saved window) ----
$ws.Range("E11").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("O25").Select()

$wb.Application.Calculate()
